$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the source values we need BEFORE any cells in row 7 are
# --- overwritten (some of the new row-7 values are simply re-used /
# --- re-ordered content that already lives elsewhere on the sheet).
$meetupText    = $ws.Range("E2").Value()   # "type: meetup ..."
$signinText    = $ws.Range("F2").Value()   # "type: signin ..."
$courseText    = $ws.Range("D2").Value()   # "type: course ..."
$subscribeText = $ws.Range("H2").Value()   # "type: subscribe ..."
$footerText    = $ws.Range("J7").Value()   # "type: footer ..." (unchanged)

# The two "blog" teaser strings (ser:71 / ser:69) that are being moved
# from columns C/G into columns E/I of row 7.
$blog71Text = @'
type: blog
width: 2
height: 1
ser: 71
'@

$blog69Text = @'
type: blog
width: 2
height: 1
ser: 69
'@

# Brand new content introduced by this edit: a new blog teaser (article 72)
# and a refreshed "Featured Video" block.
$blog72Text = @'
type: blog
width: 2
height: 1
ser: 72
'@

$videoText = @'
type: video
width: 2
height: 1
h3: Featured Video
p: May be I am wrong. But Quran is not wrong guys. I found it literally amazing. Listen to this video by Omar Sulaiman. It will give you chills. 
button.youtube: Watch on Youtube*`https://www.youtube.com/watch?v=RGypeGiH4ys&t=623s`
'@

# --- Rewrite row 7 with the new layout -------------------------------
# (new shared-strings are appended in the order they are first written,
# so write C7 -> "ser: 72" before B7 -> the new video text to land them
# on shared-string indices 20 and 21 respectively.)
$ws.Range("C7").Value = $blog72Text
$ws.Range("B7").Value = $videoText
$ws.Range("D7").Value = $meetupText
$ws.Range("E7").Value = $blog71Text
$ws.Range("F7").Value = $signinText
$ws.Range("G7").Value = $courseText
$ws.Range("H7").Value = $subscribeText
$ws.Range("I7").Value = $blog69Text
$ws.Range("J7").Value = $footerText

# --- Row height shrinks now that the tallest cell's content is shorter
$ws.Rows.Item(7).RowHeight = 255

# --- Selection / scroll position the author ended up with after editing
$ws.Range("B7").Select()
